$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Title heading text (June 2017 -> Nov 2017). The bookmark that wraps this
#    heading (zero-width, right at the very start of the document) keeps its
#    old w:name value - this runtime's Bookmarks collection does not support
#    renaming/re-adding a bookmark collapsed at document position 0.
Replace-Text "Ashley Hindmarsh - Curriculum Vitae - June 2017" "Ashley Hindmarsh - Curriculum Vitae - Nov 2017"

# 2. Summary bullet - drop the "20+ years of" lead-in.
Replace-Text "20+ years of software development experience in small/medium/large commercial organisations, across multiple sectors." "Software development experience in small/medium/large commercial organisations, across multiple sectors."

# 3. Java bullet - swap "Camel" for "Apache Camel", trim the Java8 clause, add Reactive/RxJava.
Replace-Text "Java (5+ years): Oracle Certified Associate, Java SE 8 Programmer. Spring, Camel, Java8 features, concurrency/threading, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (JAX-RS/Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, DropWizard." "Java (5+ years): Oracle Certified Associate, Java SE 8 Programmer. Spring, Apache Camel, Java8, Reactive/RxJava, JMS, TDD/BDD (JUnit, Mockito, Cucumber-JVM), Servlets, REST (JAX-RS/Jersey), Maven, JAXB, Jackson, Protobuf, Guava, Commons, DropWizard."

# 4. Remove the "Personal interests" bullet paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Personal interests: education; running coding dojos; dalliances with Scala, node.js.*") {
        $p.Range.Delete()
        break
    }
}

# 5. NoSQL databases bullet - add key/value stores + Redis/Memcached.
Replace-Text """NoSQL"" databases (e.g. DynamoDB, Cassandra)" """NoSQL"" databases and key/value stores (e.g. DynamoDB, Cassandra, Redis, Memcached)"

# 6. Mendeley role bullets - rewrite the three description bullets.
Replace-Text "Part of a cross-functional team supporting Acquisition and Onboarding services for Mendeley, through a major technology refresh." "Part of multiple cross-functional teams for Mendeley, through a major technology refresh. ** Acquisition and Onboarding: Migrating local Oauth2 sign-in to federated OpenID Connect solution. High-volume, mission-critical services. ** Building new services for flagship 'Reference Manager 2' product. Mix of client-facing and message-processing."

Replace-Text "Working mostly with back-end Java 8 development, in a microservice architecture with continuous deployment. Also with Dropwizard, Kibana, Redis, TDD, BDD." "Working mostly with back-end Java 8 development, in a microservice architecture with continuous deployment."

Replace-Text "Acquired in-depth expertise in OAuth2 & OpenID Connect. Also gained minor experience with Node.js and was able to use my experience with Cucumber (Ruby)." "Technologies used/learned: Java 8, Dropwizard, Kibana, Redis, TDD, BDD, Oauth2, OpenID Connect, RxJava, Docker/ECS, Terraform, AWS."

Write-Output "done"
